$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.798.86"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.313.40"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.72"
$ws.Range("E5").Value = "  +18.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.69"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.54"
$ws.Range("E10").Value = "  +7.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.88"
$ws.Range("E12").Value = "  +15.10%  "
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.79"
$ws.Range("E14").Value = "  +4.17%  "
$ws.Range("D15").Value = "2.661.62"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.863"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "2.317.73"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "43.787.84"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.70"
$ws.Range("E20").Value = "  +8.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.53"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +6.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.29"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  +16.09%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.50"
$ws.Range("E25").Value = "  +6.16%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.52"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.51"
$ws.Range("E28").Value = "  +11.56%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.48"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.99"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.78"
$ws.Range("E36").Value = "  +8.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.113"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +21.60%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0360"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.245"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.48"
$ws.Range("E42").Value = "  +11.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.87"
$ws.Range("E43").Value = "  +8.58%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.42"
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.81"
$ws.Range("E46").Value = "  +11.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.03"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.467"
$ws.Range("E50").Value = "  +11.36%  "
$ws.Range("E51").Value = "  +2.93%  "

$wb.Save()